{"js": "// Update citation placeholder tokens (Ref-XXXXXX) to their new values.\n// Several paragraphs share the same old token (Ref-DJ49KL) but must be\n// replaced with DIFFERENT new tokens, so each replacement is anchored on\n// a longer, paragraph-unique snippet of surrounding text rather than a\n// bare global find/replace of the token alone.\nconst replacements = [\n  {\n    find: \"Uhl-Bien et al. (Ref-DJ49KL)\",\n    replace: \"Uhl-Bien et al. (Ref-f529966)\"\n  },\n  {\n    find: \"Thompson, Gresh, and Hurwitz (Ref-DJ49KL)\",\n    replace: \"Thompson, Gresh, and Hurwitz (Ref-s628795)\"\n  },\n  {\n    find: \"(Ref-J7X2B9)\",\n    replace: \"(Ref-u395865)\"\n  },\n  {\n    find: \"(Ref-AB1CD2)\",\n    replace: \"(Ref-u763641)\"\n  },\n  {\n    find: \"Obholzer and Miller (Ref-DJ49KL)\",\n    replace: \"Obholzer and Miller (Ref-f442043)\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== 1) {\n    throw new Error(\n      `expected exactly 1 match for \"${find}\", found ${results.items.length}`\n    );\n  }\n\n  results.items[0].insertText(replace, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update citation placeholder tokens (Ref-XXXXXX) to their new values.\n# Several paragraphs share the same old token (Ref-DJ49KL) but must be\n# replaced with DIFFERENT new tokens, so each Find/Replace is anchored on\n# a longer, paragraph-unique snippet of surrounding text rather than a\n# bare global find/replace of the token alone.\n\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\nfunction Replace-UniqueText {\n    param(\n        [string]$FindText,\n        [string]$ReplaceWith\n    )\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $ok = $find.Execute(\n        $FindText,      # FindText\n        $true,          # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        $wdFindContinue,# Wrap\n        $false,         # Format\n        $ReplaceWith,   # ReplaceWith\n        $wdReplaceOne   # Replace\n    )\n\n    if (-not $ok) {\n        throw \"Could not find text: $FindText\"\n    }\n}\n\n# NOTE: the find/replace spans below deliberately avoid swallowing any\n# apostrophe from the surrounding prose (e.g. \"organization's future\") \u2014\n# Find.Execute's ReplaceWith path runs AutoFormat \"smart quotes\" on the\n# inserted span, which would mangle a straight \"'\" into a curly RIGHT\n# SINGLE QUOTATION MARK even though that character is outside the actual\n# edit. Anchoring tightly on \"(Ref-XXXXXXX)\" (already unique except for\n# the shared Ref-DJ49KL token, disambiguated with an apostrophe-free\n# author-name prefix) sidesteps that entirely.\nReplace-UniqueText \"Uhl-Bien et al. (Ref-DJ49KL)\" \"Uhl-Bien et al. (Ref-f529966)\"\nReplace-UniqueText \"Thompson, Gresh, and Hurwitz (Ref-DJ49KL)\" \"Thompson, Gresh, and Hurwitz (Ref-s628795)\"\nReplace-UniqueText \"(Ref-J7X2B9)\" \"(Ref-u395865)\"\nReplace-UniqueText \"(Ref-AB1CD2)\" \"(Ref-u763641)\"\nReplace-UniqueText \"Obholzer and Miller (Ref-DJ49KL)\" \"Obholzer and Miller (Ref-f442043)\"\n"}
